# Scheduled market-data refresh: update computed leve profit figures
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 64
$ws.Range("H64").Value = 4999.6665
$ws.Range("J64").Value = 4999.6665
$ws.Range("L64").Value = 4999.6665
$ws.Range("N64").Value = -5495.6665
# Row 67
$ws.Range("H67").Value = 4999.6665
$ws.Range("J67").Value = 4999.6665
$ws.Range("L67").Value = 4999.6665
$ws.Range("N67").Value = -6715.6665
# Row 74
$ws.Range("H74").Value = 7845.3335
$ws.Range("J74").Value = 8058.5454
$ws.Range("L74").Value = 8058.5454
$ws.Range("N74").Value = -9930.545399999999
# Row 77
$ws.Range("H77").Value = 7845.3335
$ws.Range("J77").Value = 8058.5454
$ws.Range("L77").Value = 40292.727
$ws.Range("N77").Value = -49652.727
# Row 132
$ws.Range("H132").Value = 3557.7354
$ws.Range("I132").Value = 1361.3448
$ws.Range("K132").Value = 4084.0344
$ws.Range("M132").Value = -1554.0344
# Row 137
$ws.Range("H137").Value = 3101.5881
$ws.Range("I137").Value = 1719.3334
$ws.Range("K137").Value = 5158.0002
$ws.Range("M137").Value = -2608.0002
# Row 138
$ws.Range("H138").Value = 3662.5637
$ws.Range("J138").Value = 3805.739
$ws.Range("L138").Value = 11417.217
$ws.Range("N138").Value = -21697.217

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 1663.3954
$ws.Range("I32").Value = 1027.5696
$ws.Range("K32").Value = 1027.5696
$ws.Range("M32").Value = -740.5696
# Row 45
$ws.Range("H45").Value = 2375.25
$ws.Range("I45").Value = 2375.25
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 2375.25
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -1998.25
$ws.Range("N45").ClearContents()
# Row 61
$ws.Range("H61").Value = 2288.2222
$ws.Range("I61").Value = 1861.9656
$ws.Range("J61").Value = 3060.8125
$ws.Range("K61").Value = 1861.9656
$ws.Range("L61").Value = 3060.8125
$ws.Range("M61").Value = -1649.9656
$ws.Range("N61").Value = -3484.8125
# Row 74
$ws.Range("H74").Value = 2400.88
$ws.Range("I74").Value = 1313.5
$ws.Range("J74").Value = 4334
$ws.Range("K74").Value = 1313.5
$ws.Range("L74").Value = 4334
$ws.Range("M74").Value = -439.5
$ws.Range("N74").Value = -6082
# Row 77
$ws.Range("H77").Value = 2400.88
$ws.Range("I77").Value = 1313.5
$ws.Range("J77").Value = 4334
$ws.Range("K77").Value = 6567.5
$ws.Range("L77").Value = 21670
$ws.Range("M77").Value = -2199.5
$ws.Range("N77").Value = -30406
# Row 110
$ws.Range("H110").Value = 1150.25
$ws.Range("I110").Value = 1163.6
$ws.Range("K110").Value = 1163.6
$ws.Range("M110").Value = 881.4000000000001
# Row 132
$ws.Range("H132").Value = 2034.5892
$ws.Range("I132").Value = 1719
$ws.Range("K132").Value = 5157
$ws.Range("M132").Value = -2627
# Row 136
$ws.Range("H136").Value = 2288.2222
$ws.Range("I136").Value = 1861.9656
$ws.Range("J136").Value = 3060.8125
$ws.Range("K136").Value = 5585.8968
$ws.Range("L136").Value = 9182.4375
$ws.Range("M136").Value = -3035.8968
$ws.Range("N136").Value = -14282.4375

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 2421.25
$ws.Range("I20").Value = 1892.5
$ws.Range("K20").Value = 1892.5
$ws.Range("M20").Value = -1645.5
# Row 134
$ws.Range("H134").Value = 1672.07
$ws.Range("I134").Value = 1175.974
$ws.Range("J134").Value = 3332.913
$ws.Range("K134").Value = 3527.922
$ws.Range("L134").Value = 9998.739
$ws.Range("M134").Value = -992.9219999999996
$ws.Range("N134").Value = -15068.739

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 4721.375
$ws.Range("I31").Value = 1901.4
$ws.Range("K31").Value = 1901.4
$ws.Range("M31").Value = -1606.4
# Row 34
$ws.Range("H34").Value = 4721.375
$ws.Range("I34").Value = 1901.4
$ws.Range("K34").Value = 1901.4
$ws.Range("M34").Value = -1699.4
# Row 132
$ws.Range("H132").Value = 386664.78
$ws.Range("I132").Value = 2155.75
$ws.Range("J132").Value = 1001879.2
$ws.Range("K132").Value = 6467.25
$ws.Range("L132").Value = 3005637.6
$ws.Range("M132").Value = -3937.25
$ws.Range("N132").Value = -3010697.6
# Row 134
$ws.Range("H134").Value = 8979.955
$ws.Range("I134").Value = 9302.174999999999
$ws.Range("K134").Value = 27906.525
$ws.Range("M134").Value = -25371.525

$ws = $wb.Worksheets.Item("CUL")
# Row 57
$ws.Range("H57").Value = 1488.3334
$ws.Range("I57").Value = 1488.3334
$ws.Range("K57").Value = 4465.0002
$ws.Range("M57").Value = -3906.0002

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 241542.95
$ws.Range("I80").Value = 456967
$ws.Range("J80").Value = 4576.5
$ws.Range("K80").Value = 456967
$ws.Range("L80").Value = 4576.5
$ws.Range("M80").Value = -455969
$ws.Range("N80").Value = -6572.5
# Row 83
$ws.Range("H83").Value = 241542.95
$ws.Range("I83").Value = 456967
$ws.Range("J83").Value = 4576.5
$ws.Range("K83").Value = 2284835
$ws.Range("L83").Value = 22882.5
$ws.Range("M83").Value = -2279843
$ws.Range("N83").Value = -32866.5

$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 1300.091
$ws.Range("I46").Value = 900
$ws.Range("J46").Value = 1528.7142
$ws.Range("K46").Value = 900
$ws.Range("L46").Value = 1528.7142
$ws.Range("M46").Value = -712
$ws.Range("N46").Value = -1904.7142
# Row 93
$ws.Range("H93").Value = 2975
$ws.Range("I93").Value = 2900
$ws.Range("K93").Value = 2900
$ws.Range("M93").Value = -1652
# Row 110
$ws.Range("H110").Value = 76265.71000000001
$ws.Range("J110").Value = 76265.71000000001
$ws.Range("L110").Value = 76265.71000000001
$ws.Range("N110").Value = -84445.71000000001
# Row 122
$ws.Range("H122").Value = 5041.2905
$ws.Range("I122").Value = 4149
$ws.Range("J122").Value = 7222.4443
$ws.Range("K122").Value = 12447
$ws.Range("L122").Value = 21667.3329
$ws.Range("M122").Value = -9997
$ws.Range("N122").Value = -26567.3329
# Row 132
$ws.Range("H132").Value = 2026.5676
$ws.Range("I132").Value = 1896.4849
$ws.Range("K132").Value = 5689.4547
$ws.Range("M132").Value = -3159.4547

$ws = $wb.Worksheets.Item("WVR")
# Row 3
$ws.Range("H3").Value = 45000
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").ClearContents()
# Row 132
$ws.Range("H132").Value = 1693.6097
$ws.Range("I132").Value = 1064.697
$ws.Range("J132").Value = 4287.875
$ws.Range("K132").Value = 3194.090999999999
$ws.Range("L132").Value = 12863.625
$ws.Range("M132").Value = -664.0909999999994
$ws.Range("N132").Value = -17923.625
# Row 136
$ws.Range("H136").Value = 3349.8445
$ws.Range("I136").Value = 2813.923
$ws.Range("K136").Value = 8441.769
$ws.Range("M136").Value = -5891.769
